$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry corresponds to a single cryptocurrency price/volume(1h) update
# (plus a couple of row re-ranking swaps around rows 42-45) from the latest
# GitHub Actions scrape of coinranking.com.
$updates = @(
    @{ Cell = 'D2'; Value = '29.676.31' },
    @{ Cell = 'E2'; Value = '  +0.75%  ' },
    @{ Cell = 'D3'; Value = '1.617.78' },
    @{ Cell = 'E3'; Value = '  +0.87%  ' },
    @{ Cell = 'D4'; Value = '0.993' },
    @{ Cell = 'E4'; Value = '  -0.56%  ' },
    @{ Cell = 'D5'; Value = '212.90' },
    @{ Cell = 'E5'; Value = '  +0.13%  ' },
    @{ Cell = 'D6'; Value = '0.520' },
    @{ Cell = 'E6'; Value = '  -0.44%  ' },
    @{ Cell = 'D7'; Value = '0.992' },
    @{ Cell = 'E7'; Value = '  -0.61%  ' },
    @{ Cell = 'D8'; Value = '28.81' },
    @{ Cell = 'E8'; Value = '  +7.57%  ' },
    @{ Cell = 'E9'; Value = '  +3.04%  ' },
    @{ Cell = 'D10'; Value = '0.0608' },
    @{ Cell = 'E10'; Value = '  +1.65%  ' },
    @{ Cell = 'D11'; Value = '0.0910' },
    @{ Cell = 'E11'; Value = '  -0.21%  ' },
    @{ Cell = 'D12'; Value = '1.848.91' },
    @{ Cell = 'E12'; Value = '  +0.83%  ' },
    @{ Cell = 'D13'; Value = '1.633.75' },
    @{ Cell = 'E13'; Value = '  +1.77%  ' },
    @{ Cell = 'D14'; Value = '0.567' },
    @{ Cell = 'E14'; Value = '  +6.00%  ' },
    @{ Cell = 'E15'; Value = '  +3.83%  ' },
    @{ Cell = 'D16'; Value = '29.706.35' },
    @{ Cell = 'E16'; Value = '  +0.68%  ' },
    @{ Cell = 'D17'; Value = '8.97' },
    @{ Cell = 'E17'; Value = '  +17.58%  ' },
    @{ Cell = 'D18'; Value = '64.46' },
    @{ Cell = 'E18'; Value = '  +1.68%  ' },
    @{ Cell = 'D19'; Value = '241.62' },
    @{ Cell = 'E19'; Value = '  -0.53%  ' },
    @{ Cell = 'D20'; Value = '0.0₃0709' },
    @{ Cell = 'E20'; Value = '  +2.69%  ' },
    @{ Cell = 'D21'; Value = '0.994' },
    @{ Cell = 'E21'; Value = '  -0.43%  ' },
    @{ Cell = 'E22'; Value = '  +3.11%  ' },
    @{ Cell = 'D23'; Value = '9.64' },
    @{ Cell = 'E23'; Value = '  +5.30%  ' },
    @{ Cell = 'D24'; Value = '2.10' },
    @{ Cell = 'E24'; Value = '  +0.72%  ' },
    @{ Cell = 'D25'; Value = '156.43' },
    @{ Cell = 'E25'; Value = '  +1.21%  ' },
    @{ Cell = 'E26'; Value = '  +2.35%  ' },
    @{ Cell = 'E27'; Value = '  +1.28%  ' },
    @{ Cell = 'D28'; Value = '6.58' },
    @{ Cell = 'E28'; Value = '  +3.04%  ' },
    @{ Cell = 'D29'; Value = '0.994' },
    @{ Cell = 'E29'; Value = '  -0.51%  ' },
    @{ Cell = 'E30'; Value = '  +1.70%  ' },
    @{ Cell = 'D31'; Value = '3.31' },
    @{ Cell = 'E31'; Value = '  +2.69%  ' },
    @{ Cell = 'D32'; Value = '1.08' },
    @{ Cell = 'E32'; Value = '  +1.59%  ' },
    @{ Cell = 'D33'; Value = '3.20' },
    @{ Cell = 'E33'; Value = '  +3.17%  ' },
    @{ Cell = 'D34'; Value = '1.436.31' },
    @{ Cell = 'E34'; Value = '  +1.29%  ' },
    @{ Cell = 'D35'; Value = '1.63' },
    @{ Cell = 'E35'; Value = '  +6.99%  ' },
    @{ Cell = 'E36'; Value = '  +2.11%  ' },
    @{ Cell = 'D37'; Value = '2.89' },
    @{ Cell = 'E37'; Value = '  +3.65%  ' },
    @{ Cell = 'E38'; Value = '  -0.99%  ' },
    @{ Cell = 'E39'; Value = '  +2.96%  ' },
    @{ Cell = 'D40'; Value = '0.556' },
    @{ Cell = 'E40'; Value = '  +3.81%  ' },
    @{ Cell = 'D41'; Value = '0.0502' },
    @{ Cell = 'E41'; Value = '  +3.72%  ' },
    @{ Cell = 'B42'; Value = 'ARBITRUM' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' },
    @{ Cell = 'D42'; Value = '0.827' },
    @{ Cell = 'E42'; Value = '  +4.33%  ' },
    @{ Cell = 'B43'; Value = 'RenderToken' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Cell = 'D43'; Value = '1.96' },
    @{ Cell = 'E43'; Value = '  +0.38%  ' },
    @{ Cell = 'B44'; Value = 'Aave' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Cell = 'D44'; Value = '69.78' },
    @{ Cell = 'E44'; Value = '  +6.49%  ' },
    @{ Cell = 'B45'; Value = 'BitcoinSV' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv' },
    @{ Cell = 'D45'; Value = '53.63' },
    @{ Cell = 'E45'; Value = '  +0.78%  ' },
    @{ Cell = 'D46'; Value = '0.993' },
    @{ Cell = 'E46'; Value = '  -0.55%  ' },
    @{ Cell = 'E47'; Value = '  +20.52%  ' },
    @{ Cell = 'D48'; Value = '5.45' },
    @{ Cell = 'E48'; Value = '  +3.06%  ' },
    @{ Cell = 'D49'; Value = '1.759.01' },
    @{ Cell = 'E49'; Value = '  +0.78%  ' },
    @{ Cell = 'D50'; Value = '87.63' },
    @{ Cell = 'E50'; Value = '  +1.20%  ' },
    @{ Cell = 'E51'; Value = '  -0.79%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "0.993") are not
    # silently reinterpreted as numbers by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
